$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the text summary in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.8 = 6630.2 pesos`n✅ 6630.2 pesos = 1.78 = 951.44 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 555.777
$ws2.Range("O10").Value = 3684.91
$ws2.Range("N12").Value = 3715
$ws2.Range("O12").Value = 533.105
